# Update crypto price & volume figures (D/E columns) for rows 2-51.
# D-column price strings are forced to text via a temporary "@" number
# format (otherwise Excel auto-converts lookalike numerics, dropping
# trailing zeros / switching to scientific notation); the style is then
# reset to "Normal" so no stray cell style/format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.543.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.728.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4802"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06179"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.729.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6090"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.531"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.550.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006958"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.952.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.515"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.788"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("E24").Value = "  -2.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("E26").Value = "  -1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.774"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.405"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.965"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.682"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04503"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6305"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9071"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.045"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.399"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.90%  "

$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.475"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3889"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.062"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.34%  "

$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.833"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.248"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3406"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.91%  "
